$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.553.95"
Set-TextValue $ws.Range("E2") "  +1.99%  "
Set-TextValue $ws.Range("D3") "1.571.14"
Set-TextValue $ws.Range("E3") "  +0.67%  "
Set-TextValue $ws.Range("D4") "0.991"
Set-TextValue $ws.Range("E4") "  -1.37%  "
Set-TextValue $ws.Range("D5") "211.55"
Set-TextValue $ws.Range("E5") "  +1.63%  "
Set-TextValue $ws.Range("E6") "  +0.67%  "
Set-TextValue $ws.Range("D7") "0.990"
Set-TextValue $ws.Range("E7") "  -1.49%  "
Set-TextValue $ws.Range("D8") "22.61"
Set-TextValue $ws.Range("E8") "  +2.59%  "
Set-TextValue $ws.Range("D9") "0.251"
Set-TextValue $ws.Range("E9") "  +0.95%  "
Set-TextValue $ws.Range("D10") "0.0597"
Set-TextValue $ws.Range("E10") "  -0.01%  "
Set-TextValue $ws.Range("D11") "0.0869"
Set-TextValue $ws.Range("E11") "  +1.71%  "
Set-TextValue $ws.Range("D12") "1.797.30"
Set-TextValue $ws.Range("E12") "  +0.73%  "
Set-TextValue $ws.Range("D13") "1.578.10"
Set-TextValue $ws.Range("E13") "  +3.45%  "
Set-TextValue $ws.Range("D14") "3.78"
Set-TextValue $ws.Range("E14") "  +1.52%  "
Set-TextValue $ws.Range("D15") "0.524"
Set-TextValue $ws.Range("E15") "  +0.90%  "
Set-TextValue $ws.Range("D16") "27.569.78"
Set-TextValue $ws.Range("E16") "  +2.03%  "
Set-TextValue $ws.Range("D17") "62.05"
Set-TextValue $ws.Range("E17") "  +0.32%  "
Set-TextValue $ws.Range("D18") "226.46"
Set-TextValue $ws.Range("E18") "  +4.84%  "
Set-TextValue $ws.Range("E19") "  +2.63%  "
Set-TextValue $ws.Range("D20") "0.0₃0706"
Set-TextValue $ws.Range("E20") "  -0.10%  "
Set-TextValue $ws.Range("E21") "  -1.27%  "
Set-TextValue $ws.Range("D22") "4.17"
Set-TextValue $ws.Range("E22") "  +0.90%  "
Set-TextValue $ws.Range("D23") "9.44"
Set-TextValue $ws.Range("E23") "  +2.68%  "
Set-TextValue $ws.Range("E24") "  +0.89%  "
Set-TextValue $ws.Range("D25") "150.71"
Set-TextValue $ws.Range("E25") "  -1.44%  "
Set-TextValue $ws.Range("D26") "6.66"
Set-TextValue $ws.Range("E26") "  +1.00%  "
Set-TextValue $ws.Range("D27") "15.26"
Set-TextValue $ws.Range("E27") "  +1.39%  "
Set-TextValue $ws.Range("E28") "  +1.87%  "
Set-TextValue $ws.Range("D29") "0.992"
Set-TextValue $ws.Range("E29") "  -1.29%  "
Set-TextValue $ws.Range("E30") "  +1.79%  "
Set-TextValue $ws.Range("D31") "0.0472"
Set-TextValue $ws.Range("E31") "  -0.17%  "
Set-TextValue $ws.Range("D32") "3.25"
Set-TextValue $ws.Range("E32") "  +0.69%  "
Set-TextValue $ws.Range("E33") "  +0.49%  "
Set-TextValue $ws.Range("D34") "1.458.08"
Set-TextValue $ws.Range("E34") "  +2.15%  "
Set-TextValue $ws.Range("D35") "1.11"
Set-TextValue $ws.Range("E35") "  +3.50%  "
Set-TextValue $ws.Range("D36") "1.63"
Set-TextValue $ws.Range("E36") "  +1.94%  "
Set-TextValue $ws.Range("D38") "0.0167"
Set-TextValue $ws.Range("E38") "  +0.22%  "
Set-TextValue $ws.Range("E39") "  +1.79%  "
Set-TextValue $ws.Range("D40") "0.818"
Set-TextValue $ws.Range("E40") "  +1.20%  "
Set-TextValue $ws.Range("B41") "MXToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D41") "2.34"
Set-TextValue $ws.Range("E41") "  +1.31%  "
Set-TextValue $ws.Range("B42") "FraxShare"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "5.79"
Set-TextValue $ws.Range("E42") "  -1.77%  "
Set-TextValue $ws.Range("D43") "0.991"
Set-TextValue $ws.Range("E43") "  -1.35%  "
Set-TextValue $ws.Range("D44") "65.55"
Set-TextValue $ws.Range("E44") "  +1.49%  "
Set-TextValue $ws.Range("D45") "0.967"
Set-TextValue $ws.Range("E45") "  -3.20%  "
Set-TextValue $ws.Range("D46") "1.80"
Set-TextValue $ws.Range("E46") "  +3.22%  "
Set-TextValue $ws.Range("D47") "1.712.09"
Set-TextValue $ws.Range("E47") "  +0.82%  "
Set-TextValue $ws.Range("D48") "86.62"
Set-TextValue $ws.Range("E48") "  -0.17%  "
Set-TextValue $ws.Range("D49") "0.0525"
Set-TextValue $ws.Range("E49") "  +1.50%  "
Set-TextValue $ws.Range("D50") "0.0₇0955"
Set-TextValue $ws.Range("E50") "  -7.71%  "
Set-TextValue $ws.Range("D51") "0.0939"
Set-TextValue $ws.Range("E51") "  -2.10%  "
